$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-21 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-22 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("82-59=", $true, $false, $false, $false, $false, $true, 1, $false, "38+38=", 2) | Out-Null
$d.Content.Find.Execute("81-36=", $true, $false, $false, $false, $false, $true, 1, $false, "19+12=", 2) | Out-Null
$d.Content.Find.Execute("72-67=", $true, $false, $false, $false, $false, $true, 1, $false, "63-8=", 2) | Out-Null
$d.Content.Find.Execute("24+67=", $true, $false, $false, $false, $false, $true, 1, $false, "92-55=", 2) | Out-Null
$d.Content.Find.Execute("63-45=", $true, $false, $false, $false, $false, $true, 1, $false, "93-14=", 2) | Out-Null
$d.Content.Find.Execute("91-67=", $true, $false, $false, $false, $false, $true, 1, $false, "16+19=", 2) | Out-Null
$d.Content.Find.Execute("86-29=", $true, $false, $false, $false, $false, $true, 1, $false, "57+34=", 2) | Out-Null
$d.Content.Find.Execute("66+8=", $true, $false, $false, $false, $false, $true, 1, $false, "26+45=", 2) | Out-Null
$d.Content.Find.Execute("61-36=", $true, $false, $false, $false, $false, $true, 1, $false, "94-65=", 2) | Out-Null
$d.Content.Find.Execute("51-37=", $true, $false, $false, $false, $false, $true, 1, $false, "52-24=", 2) | Out-Null
$d.Content.Find.Execute("37+48=", $true, $false, $false, $false, $false, $true, 1, $false, "24-7=", 2) | Out-Null
$d.Content.Find.Execute("14+48=", $true, $false, $false, $false, $false, $true, 1, $false, "79+9=", 2) | Out-Null
$d.Content.Find.Execute("29+34=", $true, $false, $false, $false, $false, $true, 1, $false, "96-47=", 2) | Out-Null
$d.Content.Find.Execute("66+18=", $true, $false, $false, $false, $false, $true, 1, $false, "9+34=", 2) | Out-Null
$d.Content.Find.Execute("25+29=", $true, $false, $false, $false, $false, $true, 1, $false, "59+18=", 2) | Out-Null
$d.Content.Find.Execute("77+4=", $true, $false, $false, $false, $false, $true, 1, $false, "23+68=", 2) | Out-Null
$d.Content.Find.Execute("8+65=", $true, $false, $false, $false, $false, $true, 1, $false, "43-25=", 2) | Out-Null
$d.Content.Find.Execute("13-4=", $true, $false, $false, $false, $false, $true, 1, $false, "35+6=", 2) | Out-Null
$d.Content.Find.Execute("80-61=", $true, $false, $false, $false, $false, $true, 1, $false, "85-8=", 2) | Out-Null
$d.Content.Find.Execute("86-27=", $true, $false, $false, $false, $false, $true, 1, $false, "94-85=", 2) | Out-Null
$d.Content.Find.Execute("36+28=", $true, $false, $false, $false, $false, $true, 1, $false, "92-19=", 2) | Out-Null
$d.Content.Find.Execute("71-54=", $true, $false, $false, $false, $false, $true, 1, $false, "25-6=", 2) | Out-Null
$d.Content.Find.Execute("94-87=", $true, $false, $false, $false, $false, $true, 1, $false, "64+28=", 2) | Out-Null
$d.Content.Find.Execute("92-38=", $true, $false, $false, $false, $false, $true, 1, $false, "7+26=", 2) | Out-Null
$d.Content.Find.Execute("94-28=", $true, $false, $false, $false, $false, $true, 1, $false, "46+49=", 2) | Out-Null
$d.Content.Find.Execute("17+35=", $true, $false, $false, $false, $false, $true, 1, $false, "32-28=", 2) | Out-Null
$d.Content.Find.Execute("75-19=", $true, $false, $false, $false, $false, $true, 1, $false, "55+37=", 2) | Out-Null
$d.Content.Find.Execute("24+57=", $true, $false, $false, $false, $false, $true, 1, $false, "92-55=", 2) | Out-Null
$d.Content.Find.Execute("34-18=", $true, $false, $false, $false, $false, $true, 1, $false, "39+49=", 2) | Out-Null
$d.Content.Find.Execute("26+7=", $true, $false, $false, $false, $false, $true, 1, $false, "3+8=", 2) | Out-Null
$d.Content.Find.Execute("58+13=", $true, $false, $false, $false, $false, $true, 1, $false, "28+13=", 2) | Out-Null
$d.Content.Find.Execute("73-45=", $true, $false, $false, $false, $false, $true, 1, $false, "35+16=", 2) | Out-Null
$d.Content.Find.Execute("51-49=", $true, $false, $false, $false, $false, $true, 1, $false, "5+76=", 2) | Out-Null
$d.Content.Find.Execute("16+7=", $true, $false, $false, $false, $false, $true, 1, $false, "26+15=", 2) | Out-Null
$d.Content.Find.Execute("15+19=", $true, $false, $false, $false, $false, $true, 1, $false, "16+79=", 2) | Out-Null
$d.Content.Find.Execute("62+9=", $true, $false, $false, $false, $false, $true, 1, $false, "66+5=", 2) | Out-Null
$d.Content.Find.Execute("28+35=", $true, $false, $false, $false, $false, $true, 1, $false, "86-57=", 2) | Out-Null
$d.Content.Find.Execute("38+56=", $true, $false, $false, $false, $false, $true, 1, $false, "31-3=", 2) | Out-Null
$d.Content.Find.Execute("44+48=", $true, $false, $false, $false, $false, $true, 1, $false, "58+7=", 2) | Out-Null
$d.Content.Find.Execute("84-26=", $true, $false, $false, $false, $false, $true, 1, $false, "8+76=", 2) | Out-Null
$d.Content.Find.Execute("52-23=", $true, $false, $false, $false, $false, $true, 1, $false, "39+57=", 2) | Out-Null
$d.Content.Find.Execute("75-48=", $true, $false, $false, $false, $false, $true, 1, $false, "94-27=", 2) | Out-Null
$d.Content.Find.Execute("24-5=", $true, $false, $false, $false, $false, $true, 1, $false, "92-48=", 2) | Out-Null
$d.Content.Find.Execute("23-17=", $true, $false, $false, $false, $false, $true, 1, $false, "67-29=", 2) | Out-Null
$d.Content.Find.Execute("7+8=", $true, $false, $false, $false, $false, $true, 1, $false, "87-48=", 2) | Out-Null
$d.Content.Find.Execute("73-58=", $true, $false, $false, $false, $false, $true, 1, $false, "32-24=", 2) | Out-Null
$d.Content.Find.Execute("71-48=", $true, $false, $false, $false, $false, $true, 1, $false, "43+48=", 2) | Out-Null
$d.Content.Find.Execute("26+26=", $true, $false, $false, $false, $false, $true, 1, $false, "35+28=", 2) | Out-Null
$d.Content.Find.Execute("61-39=", $true, $false, $false, $false, $false, $true, 1, $false, "82-64=", 2) | Out-Null
$d.Content.Find.Execute("40-28=", $true, $false, $false, $false, $false, $true, 1, $false, "86-17=", 2) | Out-Null
$d.Content.Find.Execute("71-12=", $true, $false, $false, $false, $false, $true, 1, $false, "40-37=", 2) | Out-Null
$d.Content.Find.Execute("18+34=", $true, $false, $false, $false, $false, $true, 1, $false, "37+38=", 2) | Out-Null
$d.Content.Find.Execute("44-8=", $true, $false, $false, $false, $false, $true, 1, $false, "60-31=", 2) | Out-Null
$d.Content.Find.Execute("38+18=", $true, $false, $false, $false, $false, $true, 1, $false, "48+44=", 2) | Out-Null
$d.Content.Find.Execute("71-25=", $true, $false, $false, $false, $false, $true, 1, $false, "58+15=", 2) | Out-Null
$d.Content.Find.Execute("90-76=", $true, $false, $false, $false, $false, $true, 1, $false, "92-64=", 2) | Out-Null
$d.Content.Find.Execute("68+13=", $true, $false, $false, $false, $false, $true, 1, $false, "34+7=", 2) | Out-Null
$d.Content.Find.Execute("93-28=", $true, $false, $false, $false, $false, $true, 1, $false, "91-8=", 2) | Out-Null
$d.Content.Find.Execute("42-35=", $true, $false, $false, $false, $false, $true, 1, $false, "69+18=", 2) | Out-Null
$d.Content.Find.Execute("29+52=", $true, $false, $false, $false, $false, $true, 1, $false, "6+17=", 2) | Out-Null
$d.Content.Find.Execute("17+6=", $true, $false, $false, $false, $false, $true, 1, $false, "55+29=", 2) | Out-Null
$d.Content.Find.Execute("83-77=", $true, $false, $false, $false, $false, $true, 1, $false, "36+57=", 2) | Out-Null
$d.Content.Find.Execute("37+9=", $true, $false, $false, $false, $false, $true, 1, $false, "82-18=", 2) | Out-Null
$d.Content.Find.Execute("48+18=", $true, $false, $false, $false, $false, $true, 1, $false, "68-39=", 2) | Out-Null
$d.Content.Find.Execute("38+25=", $true, $false, $false, $false, $false, $true, 1, $false, "90-73=", 2) | Out-Null
$d.Content.Find.Execute("84-77=", $true, $false, $false, $false, $false, $true, 1, $false, "82-3=", 2) | Out-Null
$d.Content.Find.Execute("16+75=", $true, $false, $false, $false, $false, $true, 1, $false, "40-24=", 2) | Out-Null
$d.Content.Find.Execute("87-78=", $true, $false, $false, $false, $false, $true, 1, $false, "80-21=", 2) | Out-Null
$d.Content.Find.Execute("64+8=", $true, $false, $false, $false, $false, $true, 1, $false, "57+4=", 2) | Out-Null
$d.Content.Find.Execute("54+27=", $true, $false, $false, $false, $false, $true, 1, $false, "46-29=", 2) | Out-Null
$d.Content.Find.Execute("23+8=", $true, $false, $false, $false, $false, $true, 1, $false, "83-58=", 2) | Out-Null
$d.Content.Find.Execute("91-53=", $true, $false, $false, $false, $false, $true, 1, $false, "44-39=", 2) | Out-Null
$d.Content.Find.Execute("89+9=", $true, $false, $false, $false, $false, $true, 1, $false, "55-7=", 2) | Out-Null
$d.Content.Find.Execute("68+19=", $true, $false, $false, $false, $false, $true, 1, $false, "77+5=", 2) | Out-Null
$d.Content.Find.Execute("7+28=", $true, $false, $false, $false, $false, $true, 1, $false, "39+12=", 2) | Out-Null
$d.Content.Find.Execute("90-6=", $true, $false, $false, $false, $false, $true, 1, $false, "38+33=", 2) | Out-Null
$d.Content.Find.Execute("56-28=", $true, $false, $false, $false, $false, $true, 1, $false, "23-6=", 2) | Out-Null
$d.Content.Find.Execute("28+69=", $true, $false, $false, $false, $false, $true, 1, $false, "55-37=", 2) | Out-Null
$d.Content.Find.Execute("65-17=", $true, $false, $false, $false, $false, $true, 1, $false, "20-5=", 2) | Out-Null
$d.Content.Find.Execute("20-19=", $true, $false, $false, $false, $false, $true, 1, $false, "17+6=", 2) | Out-Null
$d.Content.Find.Execute("51-2=", $true, $false, $false, $false, $false, $true, 1, $false, "18+48=", 2) | Out-Null
$d.Content.Find.Execute("44-16=", $true, $false, $false, $false, $false, $true, 1, $false, "9+82=", 2) | Out-Null
$d.Content.Find.Execute("70-17=", $true, $false, $false, $false, $false, $true, 1, $false, "80-9=", 2) | Out-Null
$d.Content.Find.Execute("16+69=", $true, $false, $false, $false, $false, $true, 1, $false, "61-44=", 2) | Out-Null
$d.Content.Find.Execute("62-54=", $true, $false, $false, $false, $false, $true, 1, $false, "85+6=", 2) | Out-Null
$d.Content.Find.Execute("82-55=", $true, $false, $false, $false, $false, $true, 1, $false, "33-14=", 2) | Out-Null
$d.Content.Find.Execute("52+9=", $true, $false, $false, $false, $false, $true, 1, $false, "57+18=", 2) | Out-Null
$d.Content.Find.Execute("43+49=", $true, $false, $false, $false, $false, $true, 1, $false, "39+33=", 2) | Out-Null
$d.Content.Find.Execute("85-9=", $true, $false, $false, $false, $false, $true, 1, $false, "5+38=", 2) | Out-Null
$d.Content.Find.Execute("96-38=", $true, $false, $false, $false, $false, $true, 1, $false, "56+16=", 2) | Out-Null
$d.Content.Find.Execute("16+57=", $true, $false, $false, $false, $false, $true, 1, $false, "83-16=", 2) | Out-Null
$d.Content.Find.Execute("63-25=", $true, $false, $false, $false, $false, $true, 1, $false, "28+13=", 2) | Out-Null
$d.Content.Find.Execute("47+29=", $true, $false, $false, $false, $false, $true, 1, $false, "34-16=", 2) | Out-Null
$d.Content.Find.Execute("30-26=", $true, $false, $false, $false, $false, $true, 1, $false, "30-11=", 2) | Out-Null
$d.Content.Find.Execute("74-65=", $true, $false, $false, $false, $false, $true, 1, $false, "18+53=", 2) | Out-Null
$d.Content.Find.Execute("34-15=", $true, $false, $false, $false, $false, $true, 1, $false, "45+47=", 2) | Out-Null
$d.Content.Find.Execute("25+17=", $true, $false, $false, $false, $false, $true, 1, $false, "70-51=", 2) | Out-Null
$d.Content.Find.Execute("2+89=", $true, $false, $false, $false, $false, $true, 1, $false, "28+45=", 2) | Out-Null
$d.Content.Find.Execute("63-6=", $true, $false, $false, $false, $false, $true, 1, $false, "8+23=", 2) | Out-Null
$d.Content.Find.Execute("33-26=", $true, $false, $false, $false, $false, $true, 1, $false, "37+6=", 2) | Out-Null
